$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.590.75"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "2.940.05"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'598.93"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "'145.40"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "'6.99"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "'0.0000226"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "'33.73"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "3.418.76"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "61.532.58"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "2.933.77"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'434.21"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'13.51"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'81.88"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").Value = "'10.94"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").Value = "'11.81"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -3.72%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'6.94"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "'26.75"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "0.0₃0883"
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'5.66"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "'42.42"
$ws.Range("E41").Value = "  +7.90%  "
$ws.Range("D42").Value = "'0.284"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").Value = "'0.0348"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "2.703.24"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "'134.59"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("D46").Value = "'365.26"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D48").Value = "'23.81"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "'2.01"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  -0.93%  "
